# ---------------------------------------------------------------------------
# Applies the "Added CO2 minimisation objective, misc improvements" commit:
#   - device sheet: new "co2em" column, several value tweaks, new row (heat
#     pump on node1), model renames (gen_el->source_el, sink_gas->export_gas,
#     sink_el->export_el on export rows)
#   - parameters sheet: planning_horizon 48 -> 24, new CO2_price=20 row
#   - new "profiles" sheet with hourly curve_wind / curve_const / curve_crude
#     profile data
#   - misc sheet-view / active-tab bookkeeping
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "device" sheet: rewrite the table with the new co2em column inserted
#    and the updated values.
# ---------------------------------------------------------------------------
$device = $wb.Worksheets.Item("device")

# Clear the previous table contents (A1:M16) before rewriting - the column
# layout changes (new column inserted at I) so a clean rewrite is safest.
$device.Range("A1:N17").ClearContents()

# Header row
$device.Range("A1").Value = "node"
$device.Range("B1").Value = "name"
$device.Range("C1").Value = "include"
$device.Range("D1").Value = "external"
$device.Range("E1").Value = "Pmax"
$device.Range("F1").Value = "Pmin"
$device.Range("G1").Value = "model"
$device.Range("H1").Value = "eta"
$device.Range("I1").Value = "co2em"
$device.Range("J1").Value = "fuelA"
$device.Range("K1").Value = "fuelB"
$device.Range("L1").Value = "heat"
$device.Range("M1").Value = "naturalpressure"
$device.Range("N1").Value = "comment"

# Row 2 - gas heater
$device.Range("A2").Value = "node1"
$device.Range("B2").Value = "gas heater"
$device.Range("C2").Value = 1
$device.Range("E2").Value = 8
$device.Range("F2").Value = 0
$device.Range("G2").Value = "gasheater"
$device.Range("H2").Value = 0.5

# Row 3 - gas turbine
$device.Range("A3").Value = "node1"
$device.Range("B3").Value = "gas turbine"
$device.Range("C3").Value = 1
$device.Range("E3").Value = 50
$device.Range("F3").Value = 0
$device.Range("G3").Value = "gasturbine"
$device.Range("J3").Value = 1
$device.Range("K3").Value = 1
$device.Range("L3").Value = 0.05

# Row 4 - compressor (el)
$device.Range("A4").Value = "node2"
$device.Range("B4").Value = "compressor"
$device.Range("C4").Value = 1
$device.Range("E4").Value = 1000
$device.Range("F4").Value = 0
$device.Range("G4").Value = "compressor_el"
$device.Range("H4").Value = 0.001

# Row 5 - compressor1 (gas)
$device.Range("A5").Value = "node1"
$device.Range("B5").Value = "compressor1"
$device.Range("C5").Value = 1
$device.Range("E5").Value = 1000
$device.Range("F5").Value = 0
$device.Range("G5").Value = "compressor_gas"
$device.Range("H5").Value = 0.004

# Row 6 - separator
$device.Range("A6").Value = "node3"
$device.Range("B6").Value = "separator"
$device.Range("C6").Value = 1
$device.Range("E6").Value = 15
$device.Range("F6").Value = 15
$device.Range("G6").Value = "sink_el"

# Row 7 - electric load
$device.Range("A7").Value = "node2"
$device.Range("B7").Value = "electric load"
$device.Range("C7").Value = 1
$device.Range("E7").Value = 50
$device.Range("F7").Value = 50
$device.Range("G7").Value = "sink_el"

# Row 8 - wind turbine
$device.Range("A8").Value = "node2"
$device.Range("B8").Value = "wind turb"
$device.Range("C8").Value = 1
$device.Range("D8").Value = "curve_wind"
$device.Range("E8").Value = 90
$device.Range("F8").Value = 0
$device.Range("G8").Value = "source_el"
$device.Range("I8").Value = 0

# Row 9 - heat pump (node2)
$device.Range("A9").Value = "node2"
$device.Range("B9").Value = "heat pump"
$device.Range("C9").Value = 1
$device.Range("E9").Value = 10
$device.Range("F9").Value = 0
$device.Range("G9").Value = "heatpump"
$device.Range("H9").Value = 3

# Row 10 - heat demand (node2)
$device.Range("A10").Value = "node2"
$device.Range("B10").Value = "heat demand"
$device.Range("C10").Value = 1
$device.Range("E10").Value = 5
$device.Range("F10").Value = 5
$device.Range("G10").Value = "sink_heat"

# Row 11 - heat dump
$device.Range("A11").Value = "node1"
$device.Range("B11").Value = "heat dump"
$device.Range("C11").Value = 0
$device.Range("E11").Value = 10
$device.Range("F11").Value = 0
$device.Range("G11").Value = "sink_heat"

# Row 12 - wellhead
$device.Range("A12").Value = "node4"
$device.Range("B12").Value = "wellhead"
$device.Range("C12").Value = 1
$device.Range("D12").Value = "curve_const"
$device.Range("E12").Value = 220
$device.Range("F12").Value = 220
$device.Range("G12").Value = "source_gas"
$device.Range("M12").Value = 600
$device.Range("N12").Value = "P=12.71GW for Åsgard"

# Row 13 - gas EXPORT
$device.Range("A13").Value = "node1ex"
$device.Range("B13").Value = "gas EXPORT"
$device.Range("C13").Value = 1
$device.Range("E13").Value = 1000
$device.Range("F13").Value = 20
$device.Range("G13").Value = "export_gas"

# Row 14 - el EXPORT
$device.Range("A14").Value = "node1ex"
$device.Range("B14").Value = "el EXPORT"
$device.Range("C14").Value = 1
$device.Range("E14").Value = 1000
$device.Range("F14").Value = 0
$device.Range("G14").Value = "export_el"

# Row 15 - diesel backup
$device.Range("A15").Value = "node1"
$device.Range("B15").Value = "diesel backup"
$device.Range("C15").Value = 1
$device.Range("E15").Value = 200
$device.Range("F15").Value = 0
$device.Range("G15").Value = "source_el"
$device.Range("I15").Value = 10

# Row 16 - heat demand (node1)
$device.Range("A16").Value = "node1"
$device.Range("B16").Value = "heat demand"
$device.Range("C16").Value = 1
$device.Range("E16").Value = 4
$device.Range("F16").Value = 4
$device.Range("G16").Value = "sink_heat"

# Row 17 (new) - heat pump (node1)
$device.Range("A17").Value = "node1"
$device.Range("B17").Value = "heat pump"
$device.Range("C17").Value = 1
$device.Range("E17").Value = 5
$device.Range("F17").Value = 0
$device.Range("G17").Value = "heatpump"
$device.Range("H17").Value = 3

# ---------------------------------------------------------------------------
# 2. "parameters" sheet: update planning_horizon, add CO2_price row
# ---------------------------------------------------------------------------
$params = $wb.Worksheets.Item("parameters")
$params.Range("B2").Value = 24
$params.Range("A4").Value = "CO2_price"
$params.Range("B4").Value = 20

# ---------------------------------------------------------------------------
# 3. New "profiles" sheet with hourly curve data, placed after "parameters"
# ---------------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIdx)
$profiles = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$profiles.Name = "profiles"

$profiles.Range("A1").Value = "timestep"
$profiles.Range("B1").Value = "curve_wind"
$profiles.Range("C1").Value = "curve_const"
$profiles.Range("D1").Value = "curve_crude"

$profileRows = @(
    @(0,  0.76644727800000001, 1, 1),
    @(1,  0.77168040699999996, 1, 0.95250000000000001),
    @(2,  0.77691353600000002, 1, 0.92749999999999999),
    @(3,  0.78214666499999996, 1, 0.8869999999999999),
    @(4,  0.78737979400000002, 1, 0.88300000000000001),
    @(5,  0.79261292299999997, 1, 0.87149999999999994),
    @(6,  0.797846053,         1, 0.82750000000000001),
    @(7,  0.80806383100000001, 1, 0.72),
    @(8,  0.81828160900000002, 1, 0.72),
    @(9,  0.828499388,         1, 0.83949999999999991),
    @(10, 0.83871716600000001, 1, 0.95250000000000001),
    @(11, 0.84558434000000005, 1, 0.96500000000000008),
    @(12, 0.85197045199999999, 1, 0.95299999999999996),
    @(13, 0.80084982500000002, 1, 0.96050000000000002),
    @(14, 0.73382039200000004, 1, 0.99700000000000011),
    @(15, 0.64636604600000003, 1, 1.0525),
    @(16, 0.565669058,         1, 1.1499999999999999),
    @(17, 0.47741213300000002, 1, 1.1890000000000001),
    @(18, 0.38533609000000002, 1, 1.1519999999999999),
    @(19, 0.41062632700000001, 1, 1.046),
    @(20, 0.438582729,         1, 0.96899999999999997),
    @(21, 0.46653913200000002, 1, 0.95250000000000001),
    @(22, 0.49449553400000001, 1, 0.9375),
    @(23, 0.52245193700000003, 1, 0.97300000000000009)
)

$r = 2
foreach ($row in $profileRows) {
    $profiles.Cells.Item($r, 1).Value = $row[0]
    $profiles.Cells.Item($r, 2).Value = $row[1]
    $profiles.Cells.Item($r, 3).Value = $row[2]
    $profiles.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$profiles.Range("B10").Select()

# ---------------------------------------------------------------------------
# 4. View / selection bookkeeping matching the edited workbook
# ---------------------------------------------------------------------------
$params.Range("C4").Select()

$device.Activate()
$device.Range("I16").Select()
